$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1062.4736
$ws.Range("I28").Value = 566
$ws.Range("J28").Value = 9999
$ws.Range("K28").Value = 566
$ws.Range("L28").Value = 9999
$ws.Range("M28").Value = -81
$ws.Range("N28").Value = -10969

$ws.Range("H107").Value = 1929.95
$ws.Range("I107").Value = 1682.5385
$ws.Range("J107").Value = 2389.4285
$ws.Range("K107").Value = 1682.5385
$ws.Range("L107").Value = 2389.4285
$ws.Range("M107").Value = 237.4614999999999
$ws.Range("N107").Value = -6229.4285

$ws.Range("H111").Value = 2065
$ws.Range("I111").Value = 2247.5
$ws.Range("K111").Value = 6742.5
$ws.Range("M111").Value = -3675.5

$ws.Range("H132").Value = 9575982
$ws.Range("I132").Value = 10259768
$ws.Range("K132").Value = 30779304
$ws.Range("M132").Value = -30776774

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13788.04
$ws.Range("I32").Value = 13788.04
$ws.Range("K32").Value = 13788.04
$ws.Range("M32").Value = -13501.04

$ws.Range("H45").Value = 4928.6313
$ws.Range("I45").Value = 2920.111
$ws.Range("K45").Value = 2920.111
$ws.Range("M45").Value = -2543.111

$ws.Range("H61").Value = 2184.9539
$ws.Range("I61").Value = 1590.1578
$ws.Range("J61").Value = 6422.875
$ws.Range("K61").Value = 1590.1578
$ws.Range("L61").Value = 6422.875
$ws.Range("M61").Value = -1378.1578
$ws.Range("N61").Value = -6846.875

$ws.Range("H128").Value = 84696.8
$ws.Range("J128").Value = 84696.8
$ws.Range("L128").Value = 84696.8
$ws.Range("N128").Value = -94656.8

$ws.Range("H130").Value = 74020
$ws.Range("J130").Value = 74020
$ws.Range("L130").Value = 74020
$ws.Range("N130").Value = -84060

$ws.Range("H132").Value = 3337.4285
$ws.Range("I132").Value = 3337.4285
$ws.Range("K132").Value = 10012.2855
$ws.Range("M132").Value = -7482.2855

$ws.Range("H136").Value = 2184.9539
$ws.Range("I136").Value = 1590.1578
$ws.Range("J136").Value = 6422.875
$ws.Range("K136").Value = 4770.4734
$ws.Range("L136").Value = 19268.625
$ws.Range("M136").Value = -2220.4734
$ws.Range("N136").Value = -24368.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 21215.074
$ws.Range("J107").Value = 4272.25
$ws.Range("L107").Value = 4272.25
$ws.Range("N107").Value = -8112.25

$ws.Range("H134").Value = 2217.42
$ws.Range("I134").Value = 2101.4792
$ws.Range("K134").Value = 6304.437600000001
$ws.Range("M134").Value = -3769.437600000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("M42").ClearContents()

$ws.Range("H50").Value = 36999.668
$ws.Range("J50").Value = 36999.668
$ws.Range("L50").Value = 36999.668
$ws.Range("N50").Value = -38249.668

$ws.Range("H99").Value = 4839.7617
$ws.Range("I99").Value = 4349
$ws.Range("K99").Value = 4349
$ws.Range("M99").Value = -2851

$ws.Range("H126").Value = 4839.7617
$ws.Range("I126").Value = 4349
$ws.Range("K126").Value = 13047
$ws.Range("M126").Value = -10577

$ws.Range("H132").Value = 47534.152
$ws.Range("I132").Value = 51328.668
$ws.Range("K132").Value = 153986.004
$ws.Range("M132").Value = -151456.004

$ws.Range("H134").Value = 3900
$ws.Range("I134").Value = 3270.5
$ws.Range("J134").Value = 4459.5557
$ws.Range("K134").Value = 9811.5
$ws.Range("L134").Value = 13378.6671
$ws.Range("M134").Value = -7276.5
$ws.Range("N134").Value = -18448.6671

$ws.Range("H141").Value = 48185.89
$ws.Range("J141").Value = 48185.89
$ws.Range("L141").Value = 48185.89
$ws.Range("N141").Value = -58545.89

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 388.41666
$ws.Range("J2").Value = 160.75
$ws.Range("L2").Value = 964.5
$ws.Range("N2").Value = -1190.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 384.73334
$ws.Range("I107").Value = 229.875
$ws.Range("J107").Value = 561.7143
$ws.Range("K107").Value = 229.875
$ws.Range("L107").Value = 561.7143
$ws.Range("M107").Value = 1690.125
$ws.Range("N107").Value = -4401.7143

$ws.Range("H113").Value = 2405.1177
$ws.Range("I113").Value = 2489.4285
$ws.Range("K113").Value = 2489.4285
$ws.Range("M113").Value = -319.4285

$ws.Range("H126").Value = 3058.5715
$ws.Range("I126").Value = 3165.5
$ws.Range("J126").Value = 2791.25
$ws.Range("K126").Value = 9496.5
$ws.Range("L126").Value = 8373.75
$ws.Range("M126").Value = -7026.5
$ws.Range("N126").Value = -13313.75

$ws.Range("H132").Value = 2426.2222
$ws.Range("I132").Value = 2426.2222
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7278.6666
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -4748.6666
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3115.5
$ws.Range("I40").Value = 3022.923
$ws.Range("J40").Value = 3516.6667
$ws.Range("K40").Value = 3022.923
$ws.Range("L40").Value = 3516.6667
$ws.Range("M40").Value = -2886.923
$ws.Range("N40").Value = -3788.6667

$ws.Range("H61").Value = 1000.75
$ws.Range("I61").Value = 928.53845
$ws.Range("J61").Value = 1134.8572
$ws.Range("K61").Value = 928.53845
$ws.Range("L61").Value = 1134.8572
$ws.Range("M61").Value = -726.53845
$ws.Range("N61").Value = -1538.8572

$ws.Range("H93").Value = 3944.8667
$ws.Range("J93").Value = 6001.5
$ws.Range("L93").Value = 6001.5
$ws.Range("N93").Value = -8497.5

$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("N100").ClearContents()

$ws.Range("H113").Value = 1000.75
$ws.Range("I113").Value = 928.53845
$ws.Range("J113").Value = 1134.8572
$ws.Range("K113").Value = 928.53845
$ws.Range("L113").Value = 1134.8572
$ws.Range("M113").Value = 1241.46155
$ws.Range("N113").Value = -5474.8572

$ws.Range("H122").Value = 13302.5
$ws.Range("I122").Value = 14172.385
$ws.Range("J122").Value = 1994
$ws.Range("K122").Value = 42517.155
$ws.Range("L122").Value = 5982
$ws.Range("M122").Value = -40067.155
$ws.Range("N122").Value = -10882

$ws.Range("H136").Value = 3360.1707
$ws.Range("I136").Value = 3189
$ws.Range("J136").Value = 3523.1904
$ws.Range("K136").Value = 9567
$ws.Range("L136").Value = 10569.5712
$ws.Range("M136").Value = -7017
$ws.Range("N136").Value = -15669.5712

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1029.5862
$ws.Range("I113").Value = 681.6667
$ws.Range("J113").Value = 1598.909
$ws.Range("K113").Value = 2045.0001
$ws.Range("L113").Value = 4796.727000000001
$ws.Range("M113").Value = 124.9999
$ws.Range("N113").Value = -9136.727000000001

$ws.Range("H119").Value = 25000
$ws.Range("J119").Value = 25000
$ws.Range("L119").Value = 25000
$ws.Range("N119").Value = -34676

$ws.Range("H122").Value = 12534635
$ws.Range("J122").Value = 3633.5
$ws.Range("L122").Value = 10900.5
$ws.Range("N122").Value = -15800.5

$ws.Range("H126").Value = 179929.5
$ws.Range("I126").Value = 1529.4286
$ws.Range("K126").Value = 4588.2858
$ws.Range("M126").Value = -2118.2858

$ws.Range("H132").Value = 5017264
$ws.Range("J132").Value = 2386.2856
$ws.Range("L132").Value = 7158.8568
$ws.Range("N132").Value = -12218.8568

$ws.Range("H138").Value = 99999
$ws.Range("J138").Value = 99999
$ws.Range("L138").Value = 99999
$ws.Range("N138").Value = -110279
